$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(187).Insert()

$ws.Cells.Item(187, 1).Value = 6
$ws.Cells.Item(187, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(187, 3).Value = "Metropolitana"
$ws.Cells.Item(187, 4).Value = 44606
$ws.Cells.Item(187, 5).Value = 13
$ws.Cells.Item(187, 6).Value = "Fruta"
$ws.Cells.Item(187, 7).Value = 100101
$ws.Cells.Item(187, 8).Value = "Berries"
$ws.Cells.Item(187, 9).Value = 100101001
$ws.Cells.Item(187, 10).Value = "Arándano (blue)"
$ws.Cells.Item(187, 11).Value = "Sin especificar"
$ws.Cells.Item(187, 12).Value = "Primera"
$ws.Cells.Item(187, 13).Value = 500
$ws.Cells.Item(187, 14).Value = 3400
$ws.Cells.Item(187, 15).Value = 3400
$ws.Cells.Item(187, 16).Value = 3400
$ws.Cells.Item(187, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(187, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(187, 19).Value = 1700
$ws.Cells.Item(187, 20).Value = 2
